$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44433
$ws.Cells.Item(2, 10).Value = 160
$ws.Cells.Item(2, 11).Value = 19000
$ws.Cells.Item(2, 12).Value = 20000
$ws.Cells.Item(2, 13).Value = 19500
$ws.Cells.Item(2, 16).Value = 488
$ws.Cells.Item(3, 4).Value = 44363
$ws.Cells.Item(3, 10).Value = 160
$ws.Cells.Item(3, 11).Value = 19000
$ws.Cells.Item(3, 12).Value = 20000
$ws.Cells.Item(3, 13).Value = 19500
$ws.Cells.Item(3, 16).Value = 488
$ws.Cells.Item(4, 4).Value = 44384
$ws.Cells.Item(4, 8).Value = 'Madrigal'
$ws.Cells.Item(4, 10).Value = 80
$ws.Cells.Item(4, 11).Value = 21000
$ws.Cells.Item(4, 12).Value = 22000
$ws.Cells.Item(4, 13).Value = 21500
$ws.Cells.Item(4, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(4, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(4, 16).Value = 538
$ws.Cells.Item(4, 17).Value = 40
$ws.Cells.Item(5, 4).Value = 44384
$ws.Cells.Item(5, 9).Value = 'Segunda'
$ws.Cells.Item(5, 10).Value = 30
$ws.Cells.Item(5, 11).Value = 19000
$ws.Cells.Item(5, 12).Value = 20000
$ws.Cells.Item(5, 13).Value = 19333
$ws.Cells.Item(5, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(5, 16).Value = 387
$ws.Cells.Item(5, 17).Value = 50
$ws.Cells.Item(6, 8).Value = 'Symphony'
$ws.Cells.Item(6, 10).Value = 100
$ws.Cells.Item(6, 11).Value = 20000
$ws.Cells.Item(6, 12).Value = 21000
$ws.Cells.Item(6, 13).Value = 20400
$ws.Cells.Item(6, 16).Value = 510
$ws.Cells.Item(7, 4).Value = 44356
$ws.Cells.Item(7, 8).Value = 'Argentina(o)'
$ws.Cells.Item(7, 9).Value = 'Primera'
$ws.Cells.Item(7, 10).Value = 120
$ws.Cells.Item(7, 13).Value = 19500
$ws.Cells.Item(7, 16).Value = 390
$ws.Cells.Item(8, 4).Value = 44412
$ws.Cells.Item(8, 10).Value = 240
$ws.Cells.Item(8, 11).Value = 21000
$ws.Cells.Item(8, 12).Value = 22000
$ws.Cells.Item(8, 13).Value = 21500
$ws.Cells.Item(8, 16).Value = 538
$ws.Cells.Item(10, 4).Value = 44468
$ws.Cells.Item(10, 8).Value = 'Argentina(o)'
$ws.Cells.Item(10, 10).Value = 120
$ws.Cells.Item(10, 11).Value = 17000
$ws.Cells.Item(10, 12).Value = 18000
$ws.Cells.Item(10, 13).Value = 17500
$ws.Cells.Item(10, 16).Value = 350
$ws.Cells.Item(11, 4).Value = 44377
$ws.Cells.Item(11, 10).Value = 150
$ws.Cells.Item(11, 11).Value = 20000
$ws.Cells.Item(11, 12).Value = 21000
$ws.Cells.Item(11, 13).Value = 20333
$ws.Cells.Item(11, 16).Value = 508
$ws.Cells.Item(12, 8).Value = 'Symphony'
$ws.Cells.Item(12, 10).Value = 60
$ws.Cells.Item(12, 11).Value = 21000
$ws.Cells.Item(12, 12).Value = 22000
$ws.Cells.Item(12, 13).Value = 21500
$ws.Cells.Item(12, 16).Value = 538
$ws.Cells.Item(13, 4).Value = 44483
$ws.Cells.Item(13, 8).Value = 'Madrigal'
$ws.Cells.Item(13, 10).Value = 120
$ws.Cells.Item(13, 11).Value = 14000
$ws.Cells.Item(13, 12).Value = 15000
$ws.Cells.Item(13, 13).Value = 14500
$ws.Cells.Item(13, 16).Value = 362
$ws.Cells.Item(14, 4).Value = 44426
$ws.Cells.Item(14, 8).Value = 'Madrigal'
$ws.Cells.Item(14, 10).Value = 150
$ws.Cells.Item(14, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(14, 16).Value = 488
$ws.Cells.Item(14, 17).Value = 40
$ws.Cells.Item(15, 4).Value = 44405
$ws.Cells.Item(15, 10).Value = 200
$ws.Cells.Item(15, 11).Value = 21000
$ws.Cells.Item(15, 12).Value = 22000
$ws.Cells.Item(15, 13).Value = 21500
$ws.Cells.Item(15, 16).Value = 538
$ws.Cells.Item(16, 4).Value = 44435
$ws.Cells.Item(16, 10).Value = 160
$ws.Cells.Item(16, 11).Value = 19000
$ws.Cells.Item(16, 12).Value = 20000
$ws.Cells.Item(16, 13).Value = 19500
$ws.Cells.Item(16, 16).Value = 488
$ws.Cells.Item(17, 4).Value = 44482
$ws.Cells.Item(17, 10).Value = 200
$ws.Cells.Item(18, 4).Value = 44398
$ws.Cells.Item(18, 10).Value = 170
$ws.Cells.Item(19, 4).Value = 44370
$ws.Cells.Item(19, 8).Value = 'Argentina(o)'
$ws.Cells.Item(19, 10).Value = 140
$ws.Cells.Item(19, 11).Value = 20000
$ws.Cells.Item(19, 12).Value = 21000
$ws.Cells.Item(19, 13).Value = 20429
$ws.Cells.Item(19, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(19, 16).Value = 409
$ws.Cells.Item(19, 17).Value = 50
$ws.Cells.Item(20, 4).Value = 44370
$ws.Cells.Item(20, 8).Value = 'Madrigal'
$ws.Cells.Item(20, 10).Value = 80
$ws.Cells.Item(20, 11).Value = 22000
$ws.Cells.Item(20, 12).Value = 23000
$ws.Cells.Item(20, 13).Value = 22500
$ws.Cells.Item(20, 14).Value = '$/caja 40 unidades'
$ws.Cells.Item(20, 16).Value = 562
$ws.Cells.Item(20, 17).Value = 40
$ws.Cells.Item(21, 4).Value = 44167
$ws.Cells.Item(21, 8).Value = 'Española'
$ws.Cells.Item(21, 10).Value = 160
$ws.Cells.Item(21, 11).Value = 13000
$ws.Cells.Item(21, 12).Value = 14000
$ws.Cells.Item(21, 13).Value = 13500
$ws.Cells.Item(21, 14).Value = '$/caja 30 unidades'
$ws.Cells.Item(21, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(21, 16).Value = 450
$ws.Cells.Item(21, 17).Value = 30
$ws.Cells.Item(22, 4).Value = 44419
$ws.Cells.Item(22, 8).Value = 'Symphony'
$ws.Cells.Item(22, 10).Value = 150
$ws.Cells.Item(22, 11).Value = 21000
$ws.Cells.Item(22, 12).Value = 22000
$ws.Cells.Item(22, 13).Value = 21500
$ws.Cells.Item(22, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(22, 16).Value = 430
$ws.Cells.Item(22, 17).Value = 50
$ws.Cells.Item(23, 4).Value = 44489
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 11).Value = 13000
$ws.Cells.Item(23, 12).Value = 14000
$ws.Cells.Item(23, 13).Value = 13500
$ws.Cells.Item(23, 16).Value = 338
$ws.Cells.Item(24, 4).Value = 44160
$ws.Cells.Item(24, 11).Value = 14000
$ws.Cells.Item(24, 12).Value = 15000
$ws.Cells.Item(24, 13).Value = 14500
$ws.Cells.Item(24, 16).Value = 362
